$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.049.63"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.759.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.80"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.605"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.78%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.72%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.66"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -16.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.245.27"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.87"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.655.82"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.759.97"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.21"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.85"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "356.53"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.74"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.531"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.43"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.171"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.64%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.21"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.25"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.54"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.25"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.94"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.18%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.81"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.32"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "338.07"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.21"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.22"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.48"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.75"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0589"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0256"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.632"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.102"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "135.30"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.04"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.12%  "
